$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "UnitedHealth Group Incorporated (UNH)"
$ws.Range("B22").Value = "Healthcare"
$ws.Range("C22").Value = 0.59
$ws.Range("D22").Value = "A"
$ws.Range("E22").Value = 159.18

$ws.Range("A23").Value = "Johnson & Johnson (JNJ)"
$ws.Range("B23").Value = "Healthcare"
$ws.Range("C23").Value = 0.52
$ws.Range("D23").Value = "B"
$ws.Range("E23").Value = 33.35

$ws.Range("A24").Value = "AbbVie Inc. (ABBV)"
$ws.Range("B24").Value = "Healthcare"
$ws.Range("C24").Value = 0.62
$ws.Range("D24").Value = "B"
$ws.Range("E24").Value = 145.27000000000001

$ws.Range("A25").Value = "Merck & Co., Inc. (MRK)"
$ws.Range("B25").Value = "Healthcare"
$ws.Range("C25").Value = 0.4
$ws.Range("D25").Value = "A"
$ws.Range("E25").Value = 43.4

$ws.Range("A26").Value = "Thermo Fisher Scientific Inc. (TMO)"
$ws.Range("B26").Value = "Healthcare"
$ws.Range("C26").Value = 0.78
$ws.Range("D26").Value = "B"
$ws.Range("E26").Value = 84.88

$ws.Range("A27").Value = "Eli Lilly and Company (LLY)"
$ws.Range("B27").Value = "Healthcare"
$ws.Range("C27").Value = 0.42
$ws.Range("D27").Value = "C"
$ws.Range("E27").Value = 675.16

$ws.Range("A28").Value = "Amgen Inc. (AMGN)"
$ws.Range("B28").Value = "Healthcare"
$ws.Range("C28").Value = 0.6
$ws.Range("D28").Value = "B"
$ws.Range("E28").Value = 57.09

$ws.Range("A29").Value = "Abbott Laboratories (ABT)"
$ws.Range("B29").Value = "Healthcare"
$ws.Range("C29").Value = 0.72
$ws.Range("D29").Value = "B"
$ws.Range("E29").Value = 50.64

$ws.Range("A30").Value = "Intuitive Surgical, Inc. (ISRG)"
$ws.Range("B30").Value = "Healthcare"
$ws.Range("C30").Value = 1.39
$ws.Range("D30").Value = "C"
$ws.Range("E30").Value = 193.29

$ws.Range("A31").Value = "Pfizer Inc. (PFE)"
$ws.Range("B31").Value = "Healthcare"
$ws.Range("C31").Value = 0.62
$ws.Range("D31").Value = "B"
$ws.Range("E31").Value = -7.22

$ws.Range("A32").Value = "The Williams Companies"
$ws.Range("B32").Value = "Energy"
$ws.Range("C32").Value = 1.06
$ws.Range("D32").Value = "C"
$ws.Range("E32").Value = 248.45

$ws.Range("A33").Value = "Phillips 66 (PSX)"
$ws.Range("B33").Value = "Energy"
$ws.Range("C33").Value = 1.32
$ws.Range("D33").Value = "B"
$ws.Range("E33").Value = 33.340000000000003

$ws.Range("A34").Value = "EOG Resources, Inc. (EOG)"
$ws.Range("B34").Value = "Energy"
$ws.Range("C34").Value = 1.29
$ws.Range("D34").Value = "B"
$ws.Range("E34").Value = 130.32

$ws.Range("A35").Value = "Valero Energy Corporation (VLO)"
$ws.Range("B35").Value = "Energy"
$ws.Range("C35").Value = 1.37
$ws.Range("D35").Value = "B"
$ws.Range("E35").Value = 69.11

$ws.Range("A36").Value = "Marathon Petroleum Corporation (MPC)"
$ws.Range("B36").Value = "Energy"
$ws.Range("C36").Value = 1.37
$ws.Range("D36").Value = "B"
$ws.Range("E36").Value = 182.79

$ws.Range("A37").Value = "Exxon Mobil Corporation (XOM)"
$ws.Range("B37").Value = "Energy"
$ws.Range("C37").Value = 0.88
$ws.Range("D37").Value = "A"
$ws.Range("E37").Value = 119.31

$ws.Range("A38").Value = "Chevron Corporation (CVX)"
$ws.Range("B38").Value = "Energy"
$ws.Range("C38").Value = 1.08
$ws.Range("D38").Value = "B"
$ws.Range("E38").Value = 60.29

$ws.Range("A39").Value = "Kinder Morgan, Inc. (KMI)"
$ws.Range("B39").Value = "Energy"
$ws.Range("C39").Value = 0.91
$ws.Range("D39").Value = "C"
$ws.Range("E39").Value = 86

$ws.Range("A40").Value = "Schlumberger Limited (SLB)"
$ws.Range("B40").Value = "Energy"
$ws.Range("C40").Value = 1.53
$ws.Range("D40").Value = "C"
$ws.Range("E40").Value = 37.72

$ws.Range("A41").Value = "ONEOK, Inc. (OKE)"
$ws.Range("B41").Value = "Energy"
$ws.Range("C41").Value = 1.62
$ws.Range("D41").Value = "C"
$ws.Range("E41").Value = 117.73

$ws.Range("A22:A25").Select() | Out-Null
